# Update countries & provincias Spain
# Refresh the covid-19 "Pais" dashboard with a newer data pull:
#  - a handful of countries swapped rank (their whole row, name + stats,
#    trades places with the neighbouring row) because the new counts
#    pushed them past the country that used to be ranked just above them
#  - several other rows just get refreshed case/recovered/death counts
#  - the "last updated" banner timestamp moves from 12:24 to 13:41

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: swap full rows whose country ranking order changed ---
$row35 = $ws.Range("A35:H35").Value2
$row36 = $ws.Range("A36:H36").Value2
$ws.Range("A35:H35").Value2 = $row36
$ws.Range("A36:H36").Value2 = $row35

$row73 = $ws.Range("A73:H73").Value2
$row74 = $ws.Range("A74:H74").Value2
$ws.Range("A73:H73").Value2 = $row74
$ws.Range("A74:H74").Value2 = $row73

$row206 = $ws.Range("A206:H206").Value2
$row207 = $ws.Range("A207:H207").Value2
$ws.Range("A206:H206").Value2 = $row207
$ws.Range("A207:H207").Value2 = $row206

$row210 = $ws.Range("A210:H210").Value2
$row211 = $ws.Range("A211:H211").Value2
$ws.Range("A210:H210").Value2 = $row211
$ws.Range("A211:H211").Value2 = $row210

$row213 = $ws.Range("A213:H213").Value2
$row214 = $ws.Range("A214:H214").Value2
$ws.Range("A213:H213").Value2 = $row214
$ws.Range("A214:H214").Value2 = $row213

# --- Step 2: apply updated case counts (post-swap) ---
$ws.Range("B4").Value = 2208787
$ws.Range("C4").Value = 387
$ws.Range("E4").Value = 1186506
$ws.Range("G4").Value = 13
$ws.Range("H4").Value = 119145

$ws.Range("B12").Value = 195051
$ws.Range("C12").Value = 2612
$ws.Range("D12").Value = 154812
$ws.Range("E12").Value = 31054
$ws.Range("G12").Value = 120
$ws.Range("H12").Value = 9185

$ws.Range("B35").Value = 37533
$ws.Range("C35").Value = 575
$ws.Range("D35").Value = 28896
$ws.Range("E35").Value = 8331
$ws.Range("G35").Value = 3
$ws.Range("H35").Value = 306

$ws.Range("B39").Value = 31183
$ws.Range("C39").Value = 29
$ws.Range("E39").Value = 327
$ws.Range("G39").Value = 2
$ws.Range("H39").Value = 1956

$ws.Range("B54").Value = 17203
$ws.Range("C54").Value = 14
$ws.Range("D54").Value = 16099
$ws.Range("E54").Value = 417
$ws.Range("G54").Value = 6
$ws.Range("H54").Value = 687

$ws.Range("B73").Value = 7177
$ws.Range("C73").Value = 586
$ws.Range("D73").Value = 1167
$ws.Range("E73").Value = 5990
$ws.Range("G73").Value = 1
$ws.Range("H73").Value = 20

$ws.Range("B77").Value = 5369
$ws.Range("C77").Value = 122
$ws.Range("D77").Value = 3606
$ws.Range("E77").Value = 1690
$ws.Range("G77").Value = 3
$ws.Range("H77").Value = 73

$ws.Range("B120").Value = 1378
$ws.Range("C120").Value = 61
$ws.Range("D120").Value = 450
$ws.Range("E120").Value = 916

$ws.Range("B124").Value = 1121
$ws.Range("C124").Value = 8
$ws.Range("D124").Value = 1071
$ws.Range("E124").Value = 46

$ws.Range("B139").Value = 696
$ws.Range("C139").Value = 2
$ws.Range("D139").Value = 591
$ws.Range("E139").Value = 63

$ws.Range("B141").Value = 662
$ws.Range("C141").Value = 6
$ws.Range("D141").Value = 610
$ws.Range("E141").Value = 43

$ws.Range("D210").Value = 10
$ws.Range("H210").Value = 1

# --- Step 3: update the "last updated" timestamp banner ---
$ws.Range("A1").Value = "Datos actualizados a 17 de Junio de 2020 a las 13:41"
